# Auto-generated Excel COM-interop script
# Applies cell value updates to the Valefor_Profits workbook sheets
# as described by the commit diff (scheduled market-data refresh).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 35723108
$ws.Range("I33").Value = 40001840
$ws.Range("K33").Value = 40001840
$ws.Range("M33").Value = -40001611

$ws.Range("H64").Value = 3230.125
$ws.Range("I64").Value = 3031.3333
$ws.Range("J64").Value = 3485.7144
$ws.Range("K64").Value = 3031.3333
$ws.Range("L64").Value = 3485.7144
$ws.Range("M64").Value = -2783.3333
$ws.Range("N64").Value = -3981.7144

$ws.Range("H67").Value = 3230.125
$ws.Range("I67").Value = 3031.3333
$ws.Range("J67").Value = 3485.7144
$ws.Range("K67").Value = 3031.3333
$ws.Range("L67").Value = 3485.7144
$ws.Range("M67").Value = -2173.3333
$ws.Range("N67").Value = -5201.7144

$ws.Range("H70").Value = 1500
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 1500
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 4500
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -5040

$ws.Range("H73").Value = 1500
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 1500
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 4500
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -6372

$ws.Range("H96").Value = 1256.875
$ws.Range("I96").Value = 1647.625
$ws.Range("J96").Value = 866.125
$ws.Range("K96").Value = 4942.875
$ws.Range("L96").Value = 2598.375
$ws.Range("M96").Value = -3569.875
$ws.Range("N96").Value = -5344.375

$ws.Range("H101").Value = 0
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("M101").ClearContents()
$ws.Range("N101").ClearContents()

$ws.Range("H113").Value = 1114801.1
$ws.Range("I113").Value = 5557305.5
$ws.Range("J113").Value = 4175
$ws.Range("K113").Value = 5557305.5
$ws.Range("L113").Value = 4175
$ws.Range("M113").Value = -5554051.5
$ws.Range("N113").Value = -10683

$ws.Range("H125").Value = 1717.5294
$ws.Range("I125").Value = 664.2857
$ws.Range("K125").Value = 5978.571300000001
$ws.Range("M125").Value = -3518.571300000001

$ws.Range("H138").Value = 3098.77
$ws.Range("I138").Value = 929
$ws.Range("J138").Value = 4545.283
$ws.Range("K138").Value = 2787
$ws.Range("L138").Value = 13635.849
$ws.Range("M138").Value = 2353
$ws.Range("N138").Value = -23915.849

$ws.Range("H141").Value = 16774258
$ws.Range("I141").Value = 21681394
$ws.Range("J141").Value = 8208.333000000001
$ws.Range("K141").Value = 65044182
$ws.Range("L141").Value = 24624.999
$ws.Range("M141").Value = -65039002
$ws.Range("N141").Value = -34984.999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 6174065.5
$ws.Range("I74").Value = 7353975
$ws.Range("J74").Value = 2230.6924
$ws.Range("K74").Value = 7353975
$ws.Range("L74").Value = 2230.6924
$ws.Range("M74").Value = -7353101
$ws.Range("N74").Value = -3978.6924

$ws.Range("H77").Value = 6174065.5
$ws.Range("I77").Value = 7353975
$ws.Range("J77").Value = 2230.6924
$ws.Range("K77").Value = 36769875
$ws.Range("L77").Value = 11153.462
$ws.Range("M77").Value = -36765507
$ws.Range("N77").Value = -19889.462

$ws.Range("H122").Value = 10418147
$ws.Range("I122").Value = 13889996
$ws.Range("J122").Value = 2600
$ws.Range("K122").Value = 41669988
$ws.Range("L122").Value = 7800
$ws.Range("M122").Value = -41667538
$ws.Range("N122").Value = -12700

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 776.7826
$ws.Range("I107").Value = 820.3333
$ws.Range("J107").Value = 620
$ws.Range("K107").Value = 820.3333
$ws.Range("L107").Value = 620
$ws.Range("M107").Value = 1099.6667
$ws.Range("N107").Value = -4460

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 270.66666
$ws.Range("I19").Value = 243.63637
$ws.Range("J19").Value = 345
$ws.Range("K19").Value = 243.63637
$ws.Range("L19").Value = 345
$ws.Range("M19").Value = -73.63637
$ws.Range("N19").Value = -685

$ws.Range("H24").Value = 270.66666
$ws.Range("I24").Value = 243.63637
$ws.Range("J24").Value = 345
$ws.Range("K24").Value = 243.63637
$ws.Range("L24").Value = 345
$ws.Range("M24").Value = -73.63637
$ws.Range("N24").Value = -685

$ws.Range("H31").Value = 14928475
$ws.Range("I31").Value = 29412898
$ws.Range("J31").Value = 5128.9697
$ws.Range("K31").Value = 29412898
$ws.Range("L31").Value = 5128.9697
$ws.Range("M31").Value = -29412603
$ws.Range("N31").Value = -5718.9697

$ws.Range("H34").Value = 14928475
$ws.Range("I34").Value = 29412898
$ws.Range("J34").Value = 5128.9697
$ws.Range("K34").Value = 29412898
$ws.Range("L34").Value = 5128.9697
$ws.Range("M34").Value = -29412696
$ws.Range("N34").Value = -5532.9697

$ws.Range("H62").Value = 3154.1
$ws.Range("I62").Value = 2745
$ws.Range("J62").Value = 3329.4285
$ws.Range("K62").Value = 2745
$ws.Range("L62").Value = 3329.4285
$ws.Range("M62").Value = -2121
$ws.Range("N62").Value = -4577.4285

$ws.Range("H65").Value = 3154.1
$ws.Range("I65").Value = 2745
$ws.Range("J65").Value = 3329.4285
$ws.Range("K65").Value = 13725
$ws.Range("L65").Value = 16647.1425
$ws.Range("M65").Value = -10605
$ws.Range("N65").Value = -22887.1425

$ws.Range("H122").Value = 16603.428
$ws.Range("I122").Value = 25556
$ws.Range("J122").Value = 4666.6665
$ws.Range("K122").Value = 76668
$ws.Range("L122").Value = 13999.9995
$ws.Range("M122").Value = -74218
$ws.Range("N122").Value = -18899.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 41667716
$ws.Range("I17").Value = 62501076
$ws.Range("J17").Value = 999.5
$ws.Range("K17").Value = 187503228
$ws.Range("L17").Value = 2998.5
$ws.Range("M17").Value = -187503059
$ws.Range("N17").Value = -3336.5

$ws.Range("H33").Value = 130.61539
$ws.Range("I33").Value = 122.833336
$ws.Range("J33").Value = 137.28572
$ws.Range("K33").Value = 737.000016
$ws.Range("L33").Value = 823.71432
$ws.Range("M33").Value = -454.000016
$ws.Range("N33").Value = -1389.71432

$ws.Range("H70").Value = 3355.2856
$ws.Range("I70").Value = 2121.75
$ws.Range("J70").Value = 5000
$ws.Range("K70").Value = 6365.25
$ws.Range("L70").Value = 15000
$ws.Range("M70").Value = -6050.25
$ws.Range("N70").Value = -15630

$ws.Range("H73").Value = 3355.2856
$ws.Range("I73").Value = 2121.75
$ws.Range("J73").Value = 5000
$ws.Range("K73").Value = 6365.25
$ws.Range("L73").Value = 15000
$ws.Range("M73").Value = -5273.25
$ws.Range("N73").Value = -17184

$ws.Range("H75").Value = 975
$ws.Range("I75").Value = 700
$ws.Range("J75").Value = 1250
$ws.Range("K75").Value = 2100
$ws.Range("L75").Value = 3750
$ws.Range("M75").Value = -1102
$ws.Range("N75").Value = -5746

$ws.Range("H78").Value = 975
$ws.Range("I78").Value = 700
$ws.Range("J78").Value = 1250
$ws.Range("K78").Value = 6300
$ws.Range("L78").Value = 11250
$ws.Range("M78").Value = -1308
$ws.Range("N78").Value = -21234

$ws.Range("H117").Value = 3249.5715
$ws.Range("J117").Value = 3863.0908
$ws.Range("L117").Value = 11589.2724
$ws.Range("N117").Value = -18473.2724

$ws.Range("H131").Value = 10041065
$ws.Range("J131").Value = 47695.535
$ws.Range("L131").Value = 143086.605
$ws.Range("N131").Value = -153166.605

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3160.1
$ws.Range("I122").Value = 3231.5715
$ws.Range("J122").Value = 2993.3333
$ws.Range("K122").Value = 9694.7145
$ws.Range("L122").Value = 8979.999899999999
$ws.Range("M122").Value = -7244.7145
$ws.Range("N122").Value = -13879.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3115.1667
$ws.Range("I40").Value = 2276
$ws.Range("J40").Value = 3786.5
$ws.Range("K40").Value = 2276
$ws.Range("L40").Value = 3786.5
$ws.Range("M40").Value = -2140
$ws.Range("N40").Value = -4058.5

$ws.Range("H100").Value = 2103.1875
$ws.Range("I100").Value = 2104.25
$ws.Range("J100").Value = 2100
$ws.Range("K100").Value = 2104.25
$ws.Range("L100").Value = 2100
$ws.Range("M100").Value = -1563.25
$ws.Range("N100").Value = -3182

$ws.Range("H122").Value = 8320.272000000001
$ws.Range("I122").Value = 9643.883
$ws.Range("J122").Value = 3820
$ws.Range("K122").Value = 28931.649
$ws.Range("L122").Value = 11460
$ws.Range("M122").Value = -26481.649
$ws.Range("N122").Value = -16360

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1589.3889
$ws.Range("I122").Value = 1340.4
$ws.Range("J122").Value = 1900.625
$ws.Range("K122").Value = 4021.2
$ws.Range("L122").Value = 5701.875
$ws.Range("M122").Value = -1571.2
$ws.Range("N122").Value = -10601.875

$ws.Range("H131").Value = 34813
$ws.Range("J131").Value = 34813
$ws.Range("L131").Value = 34813
$ws.Range("N131").Value = -44893
